# "final re-run including fig: unmet social support"
#
# Re-run of the "Impact of the psychosis disorder on the financial
# situation" bar chart against the refreshed data pull (N=146 -> N=143).
# The chart's horizontal gridline baselines, bar tops/heights, and the
# bar/axis data-label boxes all shift vertically to match the new values;
# a handful of printed counts and percentages are updated to match.
#
# Everything in this figure lives inside one top-level group shape, so
# each target shape is reached via GroupItems.Item(<name>) and only its
# Top (and, where the bar height itself changed, Height) is touched —
# the horizontal placement (Left/Width) is untouched by this edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# --- horizontal gridlines (0/10/20/30 axis rows) ---
$sh = $grp.GroupItems.Item("pl5")
$sh.Top = 365.8563232421875

$sh = $grp.GroupItems.Item("pl6")
$sh.Top = 252.27748107910156

$sh = $grp.GroupItems.Item("pl7")
$sh.Top = 138.69866943359375

# --- bars (rc8..rc13), each repositioned/resized to its new value ---
$sh = $grp.GroupItems.Item("rc8")
$sh.Top = 423.8371887207031
$sh.Height = 55.59795379638672

$sh = $grp.GroupItems.Item("rc9")
$sh.Top = 233.21536254882812
$sh.Height = 246.2198486328125

$sh = $grp.GroupItems.Item("rc10")
$sh.Top = 352.35394287109375
$sh.Height = 127.08118438720703

$sh = $grp.GroupItems.Item("rc11")
$sh.Top = 225.27276611328125
$sh.Height = 254.16236877441406

$sh = $grp.GroupItems.Item("rc12")
$sh.Top = 201.44505310058594
$sh.Height = 277.9900817871094

$sh = $grp.GroupItems.Item("rc13")
$sh.Top = 304.6985168457031
$sh.Height = 174.7366180419922

# --- count / percentage data labels above each bar ---
$sh = $grp.GroupItems.Item("tx14")
$sh.Top = 377.96490478515625
$sh.Height = 10.051496505737305
$sh.TextFrame.TextRange.Text = "7"

$sh = $grp.GroupItems.Item("tx15")
$sh.Top = 395.1512756347656
$sh.TextFrame.TextRange.Text = "(5%)"

$sh = $grp.GroupItems.Item("tx16")
$sh.Top = 186.98875427246094
$sh.TextFrame.TextRange.Text = "31"

$sh = $grp.GroupItems.Item("tx17")
$sh.Top = 204.52944946289062
$sh.TextFrame.TextRange.Text = "(22%)"

$sh = $grp.GroupItems.Item("tx18")
$sh.Top = 306.13433837890625
$sh.TextFrame.TextRange.Text = "16"

$sh = $grp.GroupItems.Item("tx19")
$sh.Top = 323.6681213378906
$sh.TextFrame.TextRange.Text = "(11%)"

$sh = $grp.GroupItems.Item("tx20")
$sh.Top = 179.04623413085938

$sh = $grp.GroupItems.Item("tx21")
$sh.Top = 196.58685302734375

$sh = $grp.GroupItems.Item("tx22")
$sh.Top = 155.2184295654297
$sh.TextFrame.TextRange.Text = "35"

$sh = $grp.GroupItems.Item("tx23")
$sh.Top = 172.75914001464844
$sh.TextFrame.TextRange.Text = "(24%)"

$sh = $grp.GroupItems.Item("tx24")
$sh.Top = 258.65252685546875
$sh.Height = 10.225197792053223
$sh.TextFrame.TextRange.Text = "22"

$sh = $grp.GroupItems.Item("tx25")
$sh.Top = 276.0126037597656
$sh.TextFrame.TextRange.Text = "(15%)"

# --- y-axis tick labels (10/20/30) ---
$sh = $grp.GroupItems.Item("tx28")
$sh.Top = 361.3797912597656

$sh = $grp.GroupItems.Item("tx29")
$sh.Top = 247.8009490966797

$sh = $grp.GroupItems.Item("tx30")
$sh.Top = 134.2162322998047

# --- y-axis tick marks (10/20/30) ---
$sh = $grp.GroupItems.Item("pl32")
$sh.Top = 365.8563232421875

$sh = $grp.GroupItems.Item("pl33")
$sh.Top = 252.27748107910156

$sh = $grp.GroupItems.Item("pl34")
$sh.Top = 138.69866943359375

# --- subtitle sample-size annotation ---
$sh = $grp.GroupItems.Item("tx50")
$sh.TextFrame.TextRange.Text = "(N=143)"
